# "added games for March and April 2024 to SQL database and all files"
# Rename the sheet/export tab to reflect the new export date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "2024-04-04"

# The sheet ends with a blank row (161) followed by a 2-row summary block
# (162: totals, 163: averages). Clear that summary block - it gets rebuilt
# two rows further down once the new game rows are in place.
$ws.Range("A162:I163").Clear()

# New row 161: Minecraft Legends (PS5), bought 2024-04-03
$ws.Cells.Item(161, 1).Value = 940
$ws.Cells.Item(161, 2).Value = "Minecraft Legends"
$ws.Cells.Item(161, 3).Value = "PS5"
$ws.Cells.Item(161, 4).Value = 45385
$ws.Cells.Item(161, 5).Value = 39.99
$ws.Cells.Item(161, 6).Value = 0
$ws.Cells.Item(161, 7).Formula = "=E161-F161"

# New row 162: Skul: The Hero Slayer (PS4), bought 2024-04-03
$ws.Cells.Item(162, 1).Value = 941
$ws.Cells.Item(162, 2).Value = "Skul: The Hero Slayer"
$ws.Cells.Item(162, 3).Value = "PS4"
$ws.Cells.Item(162, 4).Value = 45385
$ws.Cells.Item(162, 5).Value = 16.99
$ws.Cells.Item(162, 6).Value = 0
$ws.Cells.Item(162, 7).Formula = "=E162-F162"

# Match the look of the rest of the table: centered Arial, date/currency formats.
$newRows = $ws.Range("A161:G162")
$newRows.HorizontalAlignment = -4108
$newRows.VerticalAlignment = -4107
$newRows.Font.Name = "Arial"
$ws.Range("D161:D162").NumberFormat = "yyyy-mm-dd"
$ws.Range("E161:G162").NumberFormat = "[`$€]#,##0.00"

# Rebuild the summary block two rows down (164: totals, 165: averages) over
# the now-larger data range (rows 2:162).
$ws.Cells.Item(164, 5).Formula = "=SUM(E2:E162)"
$ws.Cells.Item(164, 6).Formula = "=SUM(F2:F162)"
$ws.Cells.Item(164, 7).Formula = "=SUM(G2:G162)"
$ws.Cells.Item(164, 8).Formula = "=E164-F164"
$ws.Cells.Item(164, 9).Formula = "=COUNTA(B2:B162)"

$ws.Cells.Item(165, 5).Formula = "=E164/I164"
$ws.Cells.Item(165, 6).Formula = "=F164/I164"
$ws.Cells.Item(165, 7).Formula = "=G164/I164"
$ws.Cells.Item(165, 8).Formula = "=E165-F165"
$ws.Cells.Item(165, 9).Formula = "=I164/I164"

# Extend the AutoFilter / filtered range to cover the two new rows.
$ws.Range("A1:G162").AutoFilter()
